$wb = $excel.ActiveWorkbook

# Sheet: labor_incmon_imp_stochastic_reg (sheet3)
$ws = $wb.Worksheets.Item("labor_incmon_imp_stochastic_reg")
$ws.Range("A3").Value = 2011001.6655145618
$ws.Range("B3").Value = 251357.078125
$ws.Range("C3").Value = 457395.671875
$ws.Range("F3").Value = 4013574.75
$ws.Range("G3").Value = 2078110.4379400655
$ws.Range("H3").Value = 271465.65625
$ws.Range("J3").Value = 1055699.75

# Sheet: labor_jubpenimp_stochastic_reg (sheet4)
$ws = $wb.Worksheets.Item("labor_jubpenimp_stochastic_reg")
$ws.Range("A3").Value = 1519746.6826208543
$ws.Range("C3").Value = 251357.078125
$ws.Range("E3").Value = 2455160.75
$ws.Range("F3").Value = 2555703.75
$ws.Range("G3").Value = 1519125.402360185
$ws.Range("I3").Value = 251357.078125
$ws.Range("K3").Value = 2455160.75
$ws.Range("L3").Value = 2555703.75

# Sheet: nonlabor_imp_stochastic_reg (sheet5)
$ws = $wb.Worksheets.Item("nonlabor_imp_stochastic_reg")
$ws.Range("A3").Value = 9405235.5827523023
$ws.Range("B3").Value = 120651.3984375
$ws.Range("C3").Value = 201085.671875
$ws.Range("F3").Value = 924994.0625
$ws.Range("G3").Value = 9341860.8559432384
$ws.Range("H3").Value = 120651.3984375
$ws.Range("I3").Value = 201085.671875
$ws.Range("L3").Value = 924994.0625

# Sheet: labor_beneimp_stochastic_reg (sheet6)
$ws = $wb.Worksheets.Item("labor_beneimp_stochastic_reg")
$ws.Range("A3").Value = 868275.70093919628
$ws.Range("B3").Value = 55298.55859375
$ws.Range("C3").Value = 186483.921875
$ws.Range("D3").Value = 301628.5
$ws.Range("E3").Value = 861750.875
$ws.Range("G3").Value = 918546.60322865273
$ws.Range("H3").Value = 55298.55859375
$ws.Range("J3").Value = 312130.46875
$ws.Range("K3").Value = 861750.875
